# organizei os dados de leitura
# Normalize the "Unidade" (column B) abbreviations to a consistent
# upper-case form, removing the stray space in "FAMED/ INISA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "Faed"         = "FAED"
    "Facfan"       = "FACFAN"
    "Faodo"        = "FAODO"
    "Inqui"        = "INQUI"
    "Faeng"        = "FAENG"
    "Infi"         = "INFI"
    "Facom"        = "FACOM"
    "Famez"        = "FAMEZ"
    "FAMED/ INISA" = "FAMED/INISA"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $current = $cell.Value()
    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
